$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Default" sheet (sheet1.xml)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Default")

# Remove the old column headers (Surface, Flying Low, Flying High, Space Low,
# Space High, Recovery) that used to live in C1:G1 / B1:G1.
$ws.Range("B1:G1").ClearContents()

# A column keeps 1..10 but gets new ".5" rows inserted, so rewrite it fully.
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 4.5
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 5.5
$ws.Range("A9").Value = 6
$ws.Range("A10").Value = 6.5
$ws.Range("A11").Value = 7
$ws.Range("A12").Value = 7.5
$ws.Range("A13").Value = 8
$ws.Range("A14").Value = 8.5
$ws.Range("A15").Value = 9
$ws.Range("A16").Value = 10

# New B column classification text.
$ws.Range("B2").Value = "Homeworld"
$ws.Range("B3").Value = "Homeworld Moons"
$ws.Range("B4").Value = "Inner Planets"
$ws.Range("B5").Value = "Inner Belt and Moho"
$ws.Range("B6").Value = "Outer Belt and far-inner planets"
$ws.Range("B7").Value = "Jool"
$ws.Range("B8").Value = "Sarnus"
$ws.Range("B9").Value = "Urlum"
$ws.Range("B10").Value = "Neidon"
$ws.Range("B11").Value = "Inner Kuiper Belt"
$ws.Range("B12").Value = "Outer Kuiper Belt"
$ws.Range("B13").Value = "Inner Scattered Disc"
$ws.Range("B14").Value = "Outer Scattered Disc"
$ws.Range("B15").Value = "Deep Space"

# Note in column E.
$ws.Range("E2").Value = "Extra .5 each for high inclination or wacky orbits (comets esp)"

$ws.Range("E2").Select()

# ---------------------------------------------------------------------------
# "Stock" sheet (sheet2.xml)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Stock")

$ws2.Range("B5").Value = "Dres"
$ws2.Range("B6").Value = "Jool"
$ws2.Range("B7").Value = "Eeloo"

$ws2.Range("B8").Select()

# ---------------------------------------------------------------------------
# Switch the active sheet back to "Default"
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Select()
